$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '71.116.51'
$ws.Range("E2").Value = '  -1.48%  '

# Row 3
$ws.Range("D3").Value = '3.963.27'
$ws.Range("E3").Value = '  -1.75%  '

# Row 4
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '539.15'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.67%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '149.20'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.54%  '

# Row 7
$ws.Range("D7").Value = '3.956.48'
$ws.Range("E7").Value = '  -1.69%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.688'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.40%  '

# Row 9
$ws.Range("E9").Value = '  -0.01%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.741'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.84%  '

# Row 11
$ws.Range("E11").Value = '  -5.20%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '55.99'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +16.84%  '

# Row 13
$ws.Range("E13").Value = '  -3.21%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.66'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.67%  '

# Row 15
$ws.Range("D15").Value = '4.600.38'
$ws.Range("E15").Value = '  -1.68%  '

# Row 16
$ws.Range("D16").Value = '3.967.90'
$ws.Range("E16").Value = '  -1.93%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '20.64'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.83%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.86'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.16%  '

# Row 19
$ws.Range("E19").Value = '  -1.31%  '

# Row 20
$ws.Range("E20").Value = '  -3.92%  '

# Row 21
$ws.Range("D21").Value = '71.133.35'
$ws.Range("E21").Value = '  -1.34%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '426.86'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.87%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.59'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.52%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '97.26'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -7.53%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.25'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.21%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '14.47'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.82%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.39'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.79%  '

# Row 28
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.81'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +18.00%  '

# Row 29
$ws.Range("B29").Value = 'Filecoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.68'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.05%  '

# Row 30
$ws.Range("E30").Value = '  +1.14%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '36.54'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.16%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.80'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +15.57%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '50.58'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +19.76%  '

# Row 34
$ws.Range("B34").Value = 'Cosmos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '13.41'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.95%  '

# Row 35
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.131'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.66%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '684.02'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.88%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '65.34'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.13%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.439'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.00%  '

# Row 39
$ws.Range("D39").Value = '0.0₃0820'
$ws.Range("E39").Value = '  -4.86%  '

# Row 40
$ws.Range("E40").Value = '  -1.10%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.39'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.63%  '

# Row 42
$ws.Range("E42").Value = '  +0.09%  '

# Row 43
$ws.Range("E43").Value = '  +0.19%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0482'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.46%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.20'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.86%  '

# Row 47
$ws.Range("E47").Value = '  -0.75%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.76'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +6.08%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.34'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.61%  '

# Row 50
$ws.Range("E50").Value = '  -1.91%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000271'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.13%  '
